$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")

# New "Total" header for column T
$ws.Range("T1").Value = "Total"

# New per-row totals (column T) for the existing disease category rows (2-6)
$ws.Range("T2").Value = 87725
$ws.Range("T3").Value = 11106
$ws.Range("T4").Value = 40676
$ws.Range("T5").Value = 13647
$ws.Range("T6").Value = 55256

# New row 7: "Outros" (Others) category
$ws.Range("A7").Value = "Outros"
$outrosVals = @(6838,311,419,2004,2769,2645,3085,3533,3775,4260,4965,5325,5782,5828,5844,6640,21822,498)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "7").Value = $outrosVals[$i]
}
$ws.Range("T7").Value = 86343

# New row 8: "Total" (grand total row)
$ws.Range("A8").Value = "Total"
$totalVals = @(7665,499,652,2487,3507,3676,4800,6256,7947,10901,15488,20479,25735,28779,30203,32864,92152,663)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $totalVals[$i]
}
$ws.Range("T8").Value = 294753
